$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.481.50"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.89%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.879.46"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -1.59%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.014"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -1.69%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "316.40"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -1.42%  "

$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -1.76%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5115"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -2.05%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3947"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.22%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08423"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.34%  "

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -2.38%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "41.74"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -1.82%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.268"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.69%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.877.20"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -1.96%  "

$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -1.08%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.281"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.69%  "

$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -1.55%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001108"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.76%  "

$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.76%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06746"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -1.19%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.71"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -1.79%  "

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -1.71%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.961"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -2.62%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "28.551.16"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.85%  "

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -1.53%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.272"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.47%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.094.05"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -1.77%  "

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.88%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.73"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -1.38%  "

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -2.84%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "127.39"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.34%  "

$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -1.71%  "

$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -0.64%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.799"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -3.40%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.630"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -1.74%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.02437"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -1.61%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.06504"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -2.48%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.2191"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -1.78%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.967"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -5.75%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.271"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +0.22%  "

$ws.Range("B40").Value = "ARBITRUM"
$ws.Range("C40").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.192"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.59%  "

$ws.Range("B41").Value = "InternetComputer(DFINITY)"
$ws.Range("C41").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.081"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.73%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6448"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -2.15%  "

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.11%  "

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -1.83%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6068"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -1.84%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "13.07"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -1.20%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.717"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -1.20%  "

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.47%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.202"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -7.81%  "

$ws.Range("B50").Value = "EOS"
$ws.Range("C50").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.210"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -2.71%  "

$ws.Range("B51").Value = "Quant"
$ws.Range("C51").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "122.32"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.80%  "
